# Merge local user to new user
#
# 1) Row 4 ("Generate an intro workout"): the "Req reqs" note changes
#    from "[later] time, " to "time, " (dropping the "[later]" qualifier).
# 2) A new row is inserted after "Edit any attributes for a user" (old row 11)
#    and before "Delete a user" (old row 12), documenting a new
#    "Merge local user to new user" endpoint (PA /users/merge, needs the
#    old JWT). Everything below shifts down by one row (old row 12 becomes
#    13, ..., old row 35 becomes 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the K column text on row 4 ---
$ws.Range("K4").Value = "time, "

# --- 2) Insert the new row for the merge endpoint ---
$ws.Rows("12:12").Insert()

$ws.Range("A12").Value = $true
$ws.Range("C12").Value = "General user"
$ws.Range("D12").Value = "Merge local user to new user"
$ws.Range("E12").Value = "PA"
$ws.Range("F12").Value = "/users"
$ws.Range("G12").Value = "/merge"
$ws.Range("J12").Value = "JWT"
$ws.Range("K12").Value = "old JWT,"

# G12 is a brand new cell in a column that had no neighbouring cell to
# inherit shading from during the row insert, so explicitly pull the
# green "filled" format (reusing the existing style) from F12.
$ws.Range("F12").Copy()
$ws.Range("G12").PasteSpecial(-4122)

# Re-establish the running ID formula chain (=previous row + 1) for every
# data row, since the newly inserted blank row breaks the relative chain
# that Insert() leaves behind.
for ($r = 4; $r -le 36; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 2).Formula = "=B" + $prev + "+1"
}

# --- 3) Update the selection to mirror the saved workbook state ---
$null = $ws.Range("C37").Select()
